$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ZuMyT409"
$ws.Range("B2").Value = 23100474
$ws.Range("C2").Value = "udrovxr30"
$ws.Range("D2").Value = "Hqf!&27R"
$ws.Range("F2").Value = "tsplhwEY"
$ws.Range("G2").Value = "wKGf"
